# Natmi following Dr Hou advice
#
# The sheet originally held a single "Neutro" target-cluster row. This
# split/extends the Pomc -> Oprm1 edge stats across three target clusters:
# the existing row becomes the "M1" cluster, and two new rows are added
# for the "M2" and "Neutro" clusters (recomputed specificity figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: recompute in place for the "M1" target cluster
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3304176666666667
$ws.Range("H2").Value = 0.9912529999999999
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.138311
$ws.Range("N2").Value = 0.414933
$ws.Range("O2").Value = 0.05122185706868367
$ws.Range("P2").Value = 0.05122185706868368
$ws.Range("Q2").Value = 0.04570039789433333
$ws.Range("R2").Value = 0.4113035810489999
$ws.Range("S2").Value = 0.05122185706868367
$ws.Range("T2").Value = 0.05122185706868368

# Row 3 (new): the "M2" target cluster
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pomc"
$ws.Range("C3").Value = "Oprm1"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3304176666666667
$ws.Range("H3").Value = 0.9912529999999999
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.7121919999999999
$ws.Range("N3").Value = 2.136576
$ws.Range("O3").Value = 0.2637519563119344
$ws.Range("P3").Value = 0.2637519563119344
$ws.Range("Q3").Value = 0.2353208188586667
$ws.Range("R3").Value = 2.117887369728
$ws.Range("S3").Value = 0.2637519563119344
$ws.Range("T3").Value = 0.2637519563119344

# Row 4 (new): the "Neutro" target cluster
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pomc"
$ws.Range("C4").Value = "Oprm1"
$ws.Range("D4").Value = "Neutro"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3304176666666667
$ws.Range("H4").Value = 0.9912529999999999
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.849731
$ws.Range("N4").Value = 5.549193
$ws.Range("O4").Value = 0.6850261866193819
$ws.Range("P4").Value = 0.6850261866193819
$ws.Range("Q4").Value = 0.611183800981
$ws.Range("R4").Value = 5.500654208828999
$ws.Range("S4").Value = 0.6850261866193819
$ws.Range("T4").Value = 0.6850261866193819
